$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (originally row 26) and the "SC 92" row
# (originally row 28, which becomes row 27 once "RM 232" is removed).
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Clear cells whose values became missing in the new imputation pattern.
$ws.Range("F4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("C21").Value = ""
$ws.Range("B26").Value = ""
$ws.Range("C27").Value = ""
$ws.Range("E27").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("F30").Value = ""
$ws.Range("F32").Value = ""

# Restore/update cells whose values are now present in the new imputation pattern.
$ws.Range("C6").Value = 15.1
$ws.Range("E11").Value = -7.9
$ws.Range("F17").Value = 17.78
$ws.Range("C19").Value = 13.2
$ws.Range("C23").Value = 12.2
$ws.Range("E23").Value = -7
$ws.Range("F24").Value = 16.78
$ws.Range("E25").Value = -7.1
$ws.Range("B27").Value = -20.4
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("C29").Value = 11.2
$ws.Range("E30").Value = -5.7
$ws.Range("E33").Value = -10.7
